$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated natmi LR-pair results ("Natmi following Dr Hou advice"):
# sending/target cluster combinations expanded from 2 to 3 target
# clusters (FAPs, M1, M2) per sending cluster, and all numeric scores
# recomputed accordingly.
$arr = New-Object 'object[,]' 15,20
$arr[0,0] = "ECs"
$arr[0,1] = "Icam2"
$arr[0,2] = "Itgam"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 50.17403933333333
$arr[0,7] = 150.522118
$arr[0,8] = 0.7978131386685359
$arr[0,9] = 0.797813138668536
$arr[0,10] = 1
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.142723
$arr[0,13] = 0.428169
$arr[0,14] = 0.000470790275435748
$arr[0,15] = 0.0004707902754357479
$arr[0,16] = 7.160989415771335
$arr[0,17] = 64.448904741942
$arr[0,18] = 0.0003756026673000186
$arr[0,19] = 0.0003756026673000186
$arr[1,0] = "ECs"
$arr[1,1] = "Icam2"
$arr[1,2] = "Itgam"
$arr[1,3] = "M1"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 50.17403933333333
$arr[1,7] = 150.522118
$arr[1,8] = 0.7978131386685359
$arr[1,9] = 0.797813138668536
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 148.0881626666667
$arr[1,13] = 444.264488
$arr[1,14] = 0.4884879584272602
$arr[1,15] = 0.4884879584272601
$arr[1,16] = 7430.181298438399
$arr[1,17] = 66871.63168594558
$arr[1,18] = 0.3897221113146377
$arr[1,19] = 0.3897221113146377
$arr[2,0] = "ECs"
$arr[2,1] = "Icam2"
$arr[2,2] = "Itgam"
$arr[2,3] = "M2"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 50.17403933333333
$arr[2,7] = 150.522118
$arr[2,8] = 0.7978131386685359
$arr[2,9] = 0.797813138668536
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 154.9253336666667
$arr[2,13] = 464.776001
$arr[2,14] = 0.5110412512973043
$arr[2,15] = 0.5110412512973042
$arr[2,16] = 7773.229785121124
$arr[2,17] = 69959.06806609011
$arr[2,18] = 0.4077154246865983
$arr[2,19] = 0.4077154246865983
$arr[3,0] = "FAPs"
$arr[3,1] = "Icam2"
$arr[3,2] = "Itgam"
$arr[3,3] = "FAPs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 1.656685
$arr[3,7] = 4.970055
$arr[3,8] = 0.02634280749959451
$arr[3,9] = 0.02634280749959451
$arr[3,10] = 1
$arr[3,11] = 0.3333333333333333
$arr[3,12] = 0.142723
$arr[3,13] = 0.428169
$arr[3,14] = 0.000470790275435748
$arr[3,15] = 0.0004707902754357479
$arr[3,16] = 0.2364470532550001
$arr[3,17] = 2.128023479295
$arr[3,18] = 0.00001240193759848499
$arr[3,19] = 0.00001240193759848499
$arr[4,0] = "FAPs"
$arr[4,1] = "Icam2"
$arr[4,2] = "Itgam"
$arr[4,3] = "M1"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 1.656685
$arr[4,7] = 4.970055
$arr[4,8] = 0.02634280749959451
$arr[4,9] = 0.02634280749959451
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 148.0881626666667
$arr[4,13] = 444.264488
$arr[4,14] = 0.4884879584272602
$arr[4,15] = 0.4884879584272601
$arr[4,16] = 245.3354377674267
$arr[4,17] = 2208.01893990684
$arr[4,18] = 0.01286814425471924
$arr[4,19] = 0.01286814425471924
$arr[5,0] = "FAPs"
$arr[5,1] = "Icam2"
$arr[5,2] = "Itgam"
$arr[5,3] = "M2"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 1.656685
$arr[5,7] = 4.970055
$arr[5,8] = 0.02634280749959451
$arr[5,9] = 0.02634280749959451
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 154.9253336666667
$arr[5,13] = 464.776001
$arr[5,14] = 0.5110412512973043
$arr[5,15] = 0.5110412512973042
$arr[5,16] = 256.6624764055617
$arr[5,17] = 2309.962287650055
$arr[5,18] = 0.01346226130727679
$arr[5,19] = 0.01346226130727679
$arr[6,0] = "M1"
$arr[6,1] = "Icam2"
$arr[6,2] = "Itgam"
$arr[6,3] = "FAPs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 4.722579
$arr[6,7] = 14.167737
$arr[6,8] = 0.07509332763840293
$arr[6,9] = 0.07509332763840293
$arr[6,10] = 1
$arr[6,11] = 0.3333333333333333
$arr[6,12] = 0.142723
$arr[6,13] = 0.428169
$arr[6,14] = 0.000470790275435748
$arr[6,15] = 0.0004707902754357479
$arr[6,16] = 0.674020642617
$arr[6,17] = 6.066185783553
$arr[6,18] = 0.00003535320840227058
$arr[6,19] = 0.00003535320840227058
$arr[7,0] = "M1"
$arr[7,1] = "Icam2"
$arr[7,2] = "Itgam"
$arr[7,3] = "M1"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 4.722579
$arr[7,7] = 14.167737
$arr[7,8] = 0.07509332763840293
$arr[7,9] = 0.07509332763840293
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 148.0881626666667
$arr[7,13] = 444.264488
$arr[7,14] = 0.4884879584272602
$arr[7,15] = 0.4884879584272601
$arr[7,16] = 699.358047158184
$arr[7,17] = 6294.222424423656
$arr[7,18] = 0.03668218630959279
$arr[7,19] = 0.03668218630959279
$arr[8,0] = "M1"
$arr[8,1] = "Icam2"
$arr[8,2] = "Itgam"
$arr[8,3] = "M2"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 4.722579
$arr[8,7] = 14.167737
$arr[8,8] = 0.07509332763840293
$arr[8,9] = 0.07509332763840293
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 154.9253336666667
$arr[8,13] = 464.776001
$arr[8,14] = 0.5110412512973043
$arr[8,15] = 0.5110412512973042
$arr[8,16] = 731.6471273421929
$arr[8,17] = 6584.824146079736
$arr[8,18] = 0.03837578812040787
$arr[8,19] = 0.03837578812040787
$arr[9,0] = "M2"
$arr[9,1] = "Icam2"
$arr[9,2] = "Itgam"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 5.507876
$arr[9,7] = 16.523628
$arr[9,8] = 0.08758026854811665
$arr[9,9] = 0.08758026854811667
$arr[9,10] = 1
$arr[9,11] = 0.3333333333333333
$arr[9,12] = 0.142723
$arr[9,13] = 0.428169
$arr[9,14] = 0.000470790275435748
$arr[9,15] = 0.0004707902754357479
$arr[9,16] = 0.7861005863480002
$arr[9,17] = 7.074905277132001
$arr[9,18] = 0.00004123193875250462
$arr[9,19] = 0.00004123193875250461
$arr[10,0] = "M2"
$arr[10,1] = "Icam2"
$arr[10,2] = "Itgam"
$arr[10,3] = "M1"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 5.507876
$arr[10,7] = 16.523628
$arr[10,8] = 0.08758026854811665
$arr[10,9] = 0.08758026854811667
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 148.0881626666667
$arr[10,13] = 444.264488
$arr[10,14] = 0.4884879584272602
$arr[10,15] = 0.4884879584272601
$arr[10,16] = 815.6512370358295
$arr[10,17] = 7340.861133322465
$arr[10,18] = 0.04278190658158069
$arr[10,19] = 0.04278190658158069
$arr[11,0] = "M2"
$arr[11,1] = "Icam2"
$arr[11,2] = "Itgam"
$arr[11,3] = "M2"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 5.507876
$arr[11,7] = 16.523628
$arr[11,8] = 0.08758026854811665
$arr[11,9] = 0.08758026854811667
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 154.9253336666667
$arr[11,13] = 464.776001
$arr[11,14] = 0.5110412512973043
$arr[11,15] = 0.5110412512973042
$arr[11,16] = 853.3095270946253
$arr[11,17] = 7679.785743851628
$arr[11,18] = 0.04475713002778348
$arr[11,19] = 0.04475713002778347
$arr[12,0] = "sCs"
$arr[12,1] = "Icam2"
$arr[12,2] = "Itgam"
$arr[12,3] = "FAPs"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 0.828283
$arr[12,7] = 2.484849
$arr[12,8] = 0.01317045764534999
$arr[12,9] = 0.01317045764534999
$arr[12,10] = 1
$arr[12,11] = 0.3333333333333333
$arr[12,12] = 0.142723
$arr[12,13] = 0.428169
$arr[12,14] = 0.000470790275435748
$arr[12,15] = 0.0004707902754357479
$arr[12,16] = 0.118215034609
$arr[12,17] = 1.063935311481
$arr[12,18] = 0.000006200523382469174
$arr[12,19] = 0.000006200523382469172
$arr[13,0] = "sCs"
$arr[13,1] = "Icam2"
$arr[13,2] = "Itgam"
$arr[13,3] = "M1"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 0.828283
$arr[13,7] = 2.484849
$arr[13,8] = 0.01317045764534999
$arr[13,9] = 0.01317045764534999
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 148.0881626666667
$arr[13,13] = 444.264488
$arr[13,14] = 0.4884879584272602
$arr[13,15] = 0.4884879584272601
$arr[13,16] = 122.6589076380347
$arr[13,17] = 1103.930168742312
$arr[13,18] = 0.006433609966729715
$arr[13,19] = 0.006433609966729714
$arr[14,0] = "sCs"
$arr[14,1] = "Icam2"
$arr[14,2] = "Itgam"
$arr[14,3] = "M2"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 0.828283
$arr[14,7] = 2.484849
$arr[14,8] = 0.01317045764534999
$arr[14,9] = 0.01317045764534999
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 154.9253336666667
$arr[14,13] = 464.776001
$arr[14,14] = 0.5110412512973043
$arr[14,15] = 0.5110412512973042
$arr[14,16] = 128.3220201454277
$arr[14,17] = 1154.898181308849
$arr[14,18] = 0.006730647155237805
$arr[14,19] = 0.006730647155237804

$ws.Range("A2:T16").Value = $arr
